# "10_12: Removida duplicata de correção pra fazer"
#
# The to-do list contains two near-duplicate items about commenting the
# source code inside the "quadros" (code listings):
#   - "Comentar o código-fonte nos quadros"
#   - "colocar comentários dentro dos códigos dos quadros"
# The author removed the second (duplicate) bullet. Deleting that bullet
# in real Word also drags the (hidden) "_GoBack" last-edit bookmark to the
# new cursor position, and coalesces the two runs that used to be split by
# that bookmark in the preceding bullet ("remover" + " o" | bookmark |
# "nde está...") into a single run. We reproduce both effects here.

$d = $word.ActiveDocument

# --- locate the three bullets involved, by their (stable) text ---------
$removerParaIdx = $null
$dupParaIdx = $null
$limparParaIdx = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "remover onde está*") {
        $removerParaIdx = $i
    } elseif ($t -like "colocar comentários dentro dos códigos dos quadros*") {
        $dupParaIdx = $i
    } elseif ($t -like "limpar comentários feitos em código*") {
        $limparParaIdx = $i
    }
}

# --- 1) delete the duplicate bullet -------------------------------------
$d.Paragraphs.Item($dupParaIdx).Range.Delete()

# paragraphs after the deleted one shift up by one
if ($limparParaIdx -gt $dupParaIdx) {
    $limparParaIdx = $limparParaIdx - 1
}

# --- 2) merge the bookmark-split run in the "remover" bullet -----------
$p = $d.Paragraphs.Item($removerParaIdx)
$r = $p.Range
$afterRemover = $r.Start + 7          # length of "remover"
$tailEnd = $r.End - 1                 # exclude the paragraph mark

$tailRange = $d.Range($afterRemover, $tailEnd)
$tailText = $tailRange.Text
$tailRange.Delete()

$insertPoint = $d.Range($afterRemover, $afterRemover)
$insertPoint.InsertAfter($tailText)

$mergedRange = $d.Range($afterRemover, $afterRemover + $tailText.Length)
$mergedRange.Font.StrikeThrough = -1

# --- 3) relocate the "_GoBack" bookmark to the start of "limpar" -------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$limparStart = $d.Paragraphs.Item($limparParaIdx).Range.Start
$bmRange = $d.Range($limparStart, $limparStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
